$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some "Price" column values look like plain numbers (e.g. "158.00", "1.0000").
# Force those specific cells to Text format first so Excel keeps the exact
# string (trailing zeros, leading zeros, etc.) instead of auto-converting them
# to a numeric value.
$textCells = @("D5","D6","D7","D15","D16","D17","D19","D21","D22","D23","D24","D25","D26","D28","D29","D32","D34","D35","D36","D37","D39","D41","D42","D44","D45","D47","D48","D49","D51")
foreach ($ref in $textCells) {
    $ws.Range($ref).NumberFormat = "@"
}

# Update "Price" (D) and "Volume(1h)" (E) cell values with the latest scrape
$ws.Range("D2").Value = "29.350.31"
$ws.Range("E2").Value = "  -0.07%  "
$ws.Range("D3").Value = "1.840.65"
$ws.Range("E3").Value = "  -0.19%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "239.34"
$ws.Range("E5").Value = "  -0.44%  "
$ws.Range("D6").Value = "0.6287"
$ws.Range("E6").Value = "  -0.52%  "
$ws.Range("D7").Value = "1.0000"
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("E8").Value = "  -0.79%  "
$ws.Range("E9").Value = "  +2.29%  "
$ws.Range("E10").Value = "  -0.50%  "
$ws.Range("E11").Value = "  +0.18%  "
$ws.Range("D12").Value = "1.848.09"
$ws.Range("E12").Value = "  +0.21%  "
$ws.Range("E13").Value = "  -0.83%  "
$ws.Range("E14").Value = "  -0.71%  "
$ws.Range("D15").Value = "0.00001024"
$ws.Range("E15").Value = "  -0.52%  "
$ws.Range("D16").Value = "81.57"
$ws.Range("D17").Value = "6.236"
$ws.Range("E17").Value = "  +1.26%  "
$ws.Range("D18").Value = "29.368.90"
$ws.Range("E18").Value = "  -0.07%  "
$ws.Range("D19").Value = "229.07"
$ws.Range("E19").Value = "  -0.10%  "
$ws.Range("E20").Value = "  -0.45%  "
$ws.Range("D21").Value = "1.000"
$ws.Range("E21").Value = "  +0.14%  "
$ws.Range("D22").Value = "7.339"
$ws.Range("E22").Value = "  -1.62%  "
$ws.Range("D23").Value = "0.9999"
$ws.Range("E23").Value = "  +0.15%  "
$ws.Range("D24").Value = "158.00"
$ws.Range("E24").Value = "  -0.54%  "
$ws.Range("D25").Value = "8.465"
$ws.Range("E25").Value = "  +0.59%  "
$ws.Range("D26").Value = "0.1347"
$ws.Range("E26").Value = "  -2.34%  "
$ws.Range("E27").Value = "  -1.10%  "
$ws.Range("D28").Value = "0.07243"
$ws.Range("E28").Value = "  +13.23%  "
$ws.Range("D29").Value = "1.464"
$ws.Range("E29").Value = "  +5.65%  "
$ws.Range("E30").Value = "  +0.26%  "
$ws.Range("E31").Value = "  -1.37%  "
$ws.Range("D32").Value = "4.035"
$ws.Range("E32").Value = "  -0.26%  "
$ws.Range("E33").Value = "  +0.13%  "
$ws.Range("D34").Value = "1.139"
$ws.Range("E34").Value = "  -0.24%  "
$ws.Range("D35").Value = "0.6954"
$ws.Range("E35").Value = "  -0.56%  "
$ws.Range("D36").Value = "2.578"
$ws.Range("E36").Value = "  +0.04%  "
$ws.Range("D37").Value = "0.01838"
$ws.Range("E37").Value = "  +0.94%  "
$ws.Range("E38").Value = "  -0.94%  "
$ws.Range("D39").Value = "6.844"
$ws.Range("E39").Value = "  +3.82%  "
$ws.Range("D40").Value = "1.234.20"
$ws.Range("E40").Value = "  -1.73%  "
$ws.Range("D41").Value = "0.9297"
$ws.Range("E41").Value = "  +2.60%  "
$ws.Range("D42").Value = "0.9994"
$ws.Range("E42").Value = "  +0.12%  "
$ws.Range("D43").Value = "1.993.22"
$ws.Range("E43").Value = "  -0.63%  "
$ws.Range("D44").Value = "100.50"
$ws.Range("E44").Value = "  -0.90%  "
$ws.Range("D45").Value = "65.37"
$ws.Range("E45").Value = "  -1.44%  "
$ws.Range("E46").Value = "  +1.28%  "
$ws.Range("D47").Value = "1.705"
$ws.Range("E47").Value = "  +0.13%  "
$ws.Range("D48").Value = "6.957"
$ws.Range("E48").Value = "  -1.23%  "
$ws.Range("D49").Value = "8.917"
$ws.Range("E49").Value = "  -1.42%  "
$ws.Range("E50").Value = "  -3.88%  "
$ws.Range("D51").Value = "0.3905"
$ws.Range("E51").Value = "  -0.92%  "
